$d = $word.ActiveDocument

$replacements = @(
    @{old = "27×53="; new = "72×39="},
    @{old = "24×48="; new = "82×34="},
    @{old = "88×31="; new = "94×81="},
    @{old = "38×13="; new = "35×97="},
    @{old = "51×17="; new = "33×66="},
    @{old = "43×99="; new = "19×78="},
    @{old = "62×64="; new = "58×46="},
    @{old = "51×20="; new = "93×69="},
    @{old = "11×24="; new = "23×88="},
    @{old = "64×28="; new = "52×19="},
    @{old = "73×35="; new = "50×44="},
    @{old = "66×47="; new = "97×59="},
    @{old = "78×98="; new = "97×81="},
    @{old = "79×65="; new = "27×94="},
    @{old = "20×22="; new = "14×55="},
    @{old = "58×24="; new = "42×46="},
    @{old = "92×60="; new = "72×20="},
    @{old = "34×31="; new = "80×34="},
    @{old = "65×89="; new = "82×32="},
    @{old = "51×73="; new = "67×55="},
    @{old = "36×41="; new = "26×41="},
    @{old = "60×37="; new = "62×98="},
    @{old = "13×15="; new = "42×75="},
    @{old = "78×20="; new = "17×40="},
    @{old = "29×87="; new = "50×80="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
